$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.997.39"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.981.07"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.02"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.629"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.09"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +4.52%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +1.55%  "
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.93"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +9.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.31"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.843"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.273.07"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("E16").Value = "  +3.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.987.08"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.893.99"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.09"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.15"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.24"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.51"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("E26").Value = "  +10.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.26"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.31"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.15%  "
$ws.Range("E29").Value = "  +1.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.36"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +17.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.121"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.86"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0622"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.55"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.30"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.79%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.79"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.32"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.49"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0975"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("E41").Value = "  +1.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.18"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("E43").Value = "  +0.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.73"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.371.62"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.05"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.07%  "
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.41"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.96"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +9.75%  "
